# Scheduled market-data refresh: update Anima profit calculations
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) per sheet.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 984.5833
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 983.1818
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 2949.5454
$ws.Range("M97").Value = -2504
$ws.Range("N97").Value = -3941.5454

$ws.Range("H100").Value = 7582
$ws.Range("I100").Value = 2636.6667
$ws.Range("J100").Value = 15000
$ws.Range("K100").Value = 2636.6667
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = -2095.6667
$ws.Range("N100").Value = -16082

$ws.Range("H103").Value = 101467.4
$ws.Range("J103").Value = 1925.1428
$ws.Range("L103").Value = 5775.428400000001
$ws.Range("N103").Value = -6947.428400000001

$ws.Range("H137").Value = 1853450.8
$ws.Range("I137").Value = 3624665.8
$ws.Range("K137").Value = 10873997.4
$ws.Range("M137").Value = -10871447.4


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14496966
$ws.Range("I61").Value = 33336102
$ws.Range("J61").Value = 5323.077
$ws.Range("K61").Value = 33336102
$ws.Range("L61").Value = 5323.077
$ws.Range("M61").Value = -33335890
$ws.Range("N61").Value = -5747.077

$ws.Range("H64").Value = 98091
$ws.Range("J64").Value = 98091
$ws.Range("L64").Value = 98091
$ws.Range("N64").Value = -98587

$ws.Range("H67").Value = 98091
$ws.Range("J67").Value = 98091
$ws.Range("L67").Value = 98091
$ws.Range("N67").Value = -99807

$ws.Range("H74").Value = 10418150
$ws.Range("I74").Value = 713.71875
$ws.Range("K74").Value = 713.71875
$ws.Range("M74").Value = 160.28125

$ws.Range("H77").Value = 10418150
$ws.Range("I77").Value = 713.71875
$ws.Range("K77").Value = 3568.59375
$ws.Range("M77").Value = 799.40625

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 1352629.4
$ws.Range("I132").Value = 2268.0977
$ws.Range("J132").Value = 4812930
$ws.Range("K132").Value = 6804.293099999999
$ws.Range("L132").Value = 14438790
$ws.Range("M132").Value = -4274.293099999999
$ws.Range("N132").Value = -14443850

$ws.Range("H135").Value = 56500
$ws.Range("J135").Value = 56500
$ws.Range("L135").Value = 56500
$ws.Range("N135").Value = -66640

$ws.Range("H136").Value = 14496966
$ws.Range("I136").Value = 33336102
$ws.Range("J136").Value = 5323.077
$ws.Range("K136").Value = 100008306
$ws.Range("L136").Value = 15969.231
$ws.Range("M136").Value = -100005756
$ws.Range("N136").Value = -21069.231

$ws.Range("H137").Value = 66933
$ws.Range("J137").Value = 68036.664
$ws.Range("L137").Value = 68036.664
$ws.Range("N137").Value = -78236.664

$ws.Range("H139").Value = 84419.336
$ws.Range("J139").Value = 84419.336
$ws.Range("L139").Value = 84419.336
$ws.Range("N139").Value = -94699.336


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 61570.668
$ws.Range("J13").Value = 61570.668
$ws.Range("L13").Value = 61570.668
$ws.Range("N13").Value = -61906.668

$ws.Range("H76").Value = 20285
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 20285
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H81").Value = 54575.332
$ws.Range("J81").Value = 54575.332
$ws.Range("L81").Value = 54575.332
$ws.Range("N81").Value = -56697.332

$ws.Range("H84").Value = 54575.332
$ws.Range("J84").Value = 54575.332
$ws.Range("L84").Value = 163725.996
$ws.Range("N84").Value = -174333.996

$ws.Range("H97").Value = 16446.455
$ws.Range("I97").Value = 4937.778
$ws.Range("K97").Value = 4937.778
$ws.Range("M97").Value = -3946.778

$ws.Range("H132").Value = 77262.10000000001
$ws.Range("J132").Value = 77110
$ws.Range("L132").Value = 77110
$ws.Range("N132").Value = -87230

$ws.Range("H134").Value = 2477.6667
$ws.Range("I134").Value = 2371.0527
$ws.Range("K134").Value = 7113.158100000001
$ws.Range("M134").Value = -4578.158100000001

$ws.Range("H137").Value = 70693.336
$ws.Range("J137").Value = 70693.336
$ws.Range("L137").Value = 70693.336
$ws.Range("N137").Value = -80893.336

$ws.Range("H138").Value = 47100.285
$ws.Range("I138").Value = 5000
$ws.Range("J138").Value = 50338.77
$ws.Range("K138").Value = 5000
$ws.Range("L138").Value = 50338.77
$ws.Range("M138").Value = 140
$ws.Range("N138").Value = -60618.77


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11906410
$ws.Range("I132").Value = 12196451
$ws.Range("J132").Value = 11113633
$ws.Range("K132").Value = 36589353
$ws.Range("L132").Value = 33340899
$ws.Range("M132").Value = -36586823
$ws.Range("N132").Value = -33345959


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 5447.9473
$ws.Range("I33").Value = 7348.4287
$ws.Range("J33").Value = 126.6
$ws.Range("K33").Value = 44090.5722
$ws.Range("L33").Value = 759.5999999999999
$ws.Range("M33").Value = -43807.5722
$ws.Range("N33").Value = -1325.6

$ws.Range("H107").Value = 2250.7646
$ws.Range("I107").Value = 283.8
$ws.Range("J107").Value = 3070.3333
$ws.Range("K107").Value = 851.4000000000001
$ws.Range("L107").Value = 9210.999899999999
$ws.Range("M107").Value = 1068.6
$ws.Range("N107").Value = -13050.9999


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15125.5
$ws.Range("J92").Value = 15125.5
$ws.Range("L92").Value = 15125.5
$ws.Range("N92").Value = -18869.5

$ws.Range("H123").Value = 8565.237999999999
$ws.Range("J123").Value = 9874.706
$ws.Range("L123").Value = 9874.706
$ws.Range("N123").Value = -14774.706

$ws.Range("H132").Value = 37043212
$ws.Range("I132").Value = 52638564
$ws.Range("J132").Value = 4253.125
$ws.Range("K132").Value = 157915692
$ws.Range("L132").Value = 12759.375
$ws.Range("M132").Value = -157913162
$ws.Range("N132").Value = -17819.375


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8000
$ws.Range("J2").Value = 8000
$ws.Range("L2").Value = 8000
$ws.Range("N2").Value = -8224

$ws.Range("H132").Value = 2423.0227
$ws.Range("I132").Value = 1574.5807
$ws.Range("J132").Value = 4446.231
$ws.Range("K132").Value = 4723.742099999999
$ws.Range("L132").Value = 13338.693
$ws.Range("M132").Value = -2193.742099999999
$ws.Range("N132").Value = -18398.693

$ws.Range("H136").Value = 6175253.5
$ws.Range("I136").Value = 2043.1666
$ws.Range("K136").Value = 6129.4998
$ws.Range("M136").Value = -3579.4998


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3393254.5
$ws.Range("I132").Value = 1602.723
$ws.Range("J132").Value = 13891224
$ws.Range("K132").Value = 4808.169
$ws.Range("L132").Value = 41673672
$ws.Range("M132").Value = -2278.169
$ws.Range("N132").Value = -41678732

$ws.Range("H136").Value = 5114.759
$ws.Range("I136").Value = 5974.9165
$ws.Range("J136").Value = 4507.5884
$ws.Range("K136").Value = 17924.7495
$ws.Range("L136").Value = 13522.7652
$ws.Range("M136").Value = -15374.7495
$ws.Range("N136").Value = -18622.7652
